# Auto-generated Excel COM-interop script
# Applies scheduled-market-data refresh updates to the Pandaemonium_Profits workbook
# (columns H-N: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6500.5
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968

$ws.Range("H80").Value = 3262.9302
$ws.Range("I80").Value = 745.8125
$ws.Range("J80").Value = 4754.5557
$ws.Range("K80").Value = 2237.4375
$ws.Range("L80").Value = 14263.6671
$ws.Range("M80").Value = -1239.4375
$ws.Range("N80").Value = -16259.6671

$ws.Range("H83").Value = 3262.9302
$ws.Range("I83").Value = 745.8125
$ws.Range("J83").Value = 4754.5557
$ws.Range("K83").Value = 6712.3125
$ws.Range("L83").Value = 42791.0013
$ws.Range("M83").Value = -1720.3125
$ws.Range("N83").Value = -52775.0013

$ws.Range("H86").Value = 115555.89
$ws.Range("I86").Value = 203800.6
$ws.Range("K86").Value = 203800.6
$ws.Range("M86").Value = -202677.6

$ws.Range("H89").Value = 115555.89
$ws.Range("I89").Value = 203800.6
$ws.Range("K89").Value = 1019003
$ws.Range("M89").Value = -1013387

$ws.Range("H121").Value = 1069.4
$ws.Range("I121").Value = 1525
$ws.Range("J121").Value = 765.6667
$ws.Range("K121").Value = 4575
$ws.Range("L121").Value = 2297.0001
$ws.Range("M121").Value = -2828
$ws.Range("N121").Value = -5791.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4092.75
$ws.Range("I28").Value = 4092.75
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 4092.75
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -3900.75

$ws.Range("H32").Value = 8097.5737
$ws.Range("I32").Value = 7526.4365
$ws.Range("J32").Value = 13333
$ws.Range("K32").Value = 7526.4365
$ws.Range("L32").Value = 13333
$ws.Range("M32").Value = -7239.4365
$ws.Range("N32").Value = -13907

$ws.Range("H61").Value = 7108.6963
$ws.Range("I61").Value = 3555.4856
$ws.Range("J61").Value = 13030.714
$ws.Range("K61").Value = 3555.4856
$ws.Range("L61").Value = 13030.714
$ws.Range("M61").Value = -3343.4856
$ws.Range("N61").Value = -13454.714

$ws.Range("H88").Value = 4857.5884
$ws.Range("I88").Value = 11723.6
$ws.Range("J88").Value = 1996.75
$ws.Range("K88").Value = 11723.6
$ws.Range("L88").Value = 1996.75
$ws.Range("M88").Value = -11317.6
$ws.Range("N88").Value = -2808.75

$ws.Range("H91").Value = 4857.5884
$ws.Range("I91").Value = 11723.6
$ws.Range("J91").Value = 1996.75
$ws.Range("K91").Value = 11723.6
$ws.Range("L91").Value = 1996.75
$ws.Range("M91").Value = -10319.6
$ws.Range("N91").Value = -4804.75

$ws.Range("H99").Value = 4092.75
$ws.Range("I99").Value = 4092.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4092.75
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -1097.75

$ws.Range("H136").Value = 7108.6963
$ws.Range("I136").Value = 3555.4856
$ws.Range("J136").Value = 13030.714
$ws.Range("K136").Value = 10666.4568
$ws.Range("L136").Value = 39092.142
$ws.Range("M136").Value = -8116.4568
$ws.Range("N136").Value = -44192.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0

$ws.Range("H96").Value = 12332
$ws.Range("I96").Value = 12332
$ws.Range("K96").Value = 12332
$ws.Range("M96").Value = -9586

$ws.Range("H97").Value = 4578.364
$ws.Range("I97").Value = 4578.364
$ws.Range("K97").Value = 4578.364
$ws.Range("M97").Value = -3587.364

$ws.Range("H105").Value = 4944.44
$ws.Range("I105").Value = 4634.7827
$ws.Range("K105").Value = 4634.7827
$ws.Range("M105").Value = -2887.7827

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2089.4285
$ws.Range("J31").Value = 3416.5
$ws.Range("L31").Value = 3416.5
$ws.Range("N31").Value = -4006.5

$ws.Range("H34").Value = 2089.4285
$ws.Range("J34").Value = 3416.5
$ws.Range("L34").Value = 3416.5
$ws.Range("N34").Value = -3820.5

$ws.Range("H62").Value = 2781
$ws.Range("J62").Value = 1900
$ws.Range("L62").Value = 1900
$ws.Range("N62").Value = -3148

$ws.Range("H65").Value = 2781
$ws.Range("J65").Value = 1900
$ws.Range("L65").Value = 9500
$ws.Range("N65").Value = -15740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 734.3929000000001
$ws.Range("I113").Value = 758.7143
$ws.Range("J113").Value = 661.4286
$ws.Range("K113").Value = 2276.1429
$ws.Range("L113").Value = 1984.2858
$ws.Range("M113").Value = -106.1428999999998
$ws.Range("N113").Value = -6324.2858

$ws.Range("H140").Value = 1455.8718
$ws.Range("I140").Value = 1251.9714
$ws.Range("J140").Value = 3240
$ws.Range("K140").Value = 3755.9142
$ws.Range("L140").Value = 9720
$ws.Range("M140").Value = 1424.0858
$ws.Range("N140").Value = -20080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8627.777
$ws.Range("I80").Value = 18283.334
$ws.Range("J80").Value = 3800
$ws.Range("K80").Value = 18283.334
$ws.Range("L80").Value = 3800
$ws.Range("M80").Value = -17285.334
$ws.Range("N80").Value = -5796

$ws.Range("H83").Value = 8627.777
$ws.Range("I83").Value = 18283.334
$ws.Range("J83").Value = 3800
$ws.Range("K83").Value = 91416.67
$ws.Range("L83").Value = 19000
$ws.Range("M83").Value = -86424.67
$ws.Range("N83").Value = -28984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7950
$ws.Range("I68").Value = 7950
$ws.Range("K68").Value = 7950
$ws.Range("M68").Value = -7201

$ws.Range("H71").Value = 7950
$ws.Range("I71").Value = 7950
$ws.Range("K71").Value = 39750
$ws.Range("M71").Value = -36006

$ws.Range("H82").Value = 2216.353
$ws.Range("I82").Value = 1333.2
$ws.Range("J82").Value = 3478
$ws.Range("K82").Value = 1333.2
$ws.Range("L82").Value = 3478
$ws.Range("M82").Value = -972.2
$ws.Range("N82").Value = -4200

$ws.Range("H85").Value = 2216.353
$ws.Range("I85").Value = 1333.2
$ws.Range("J85").Value = 3478
$ws.Range("K85").Value = 1333.2
$ws.Range("L85").Value = 3478
$ws.Range("M85").Value = -85.20000000000005
$ws.Range("N85").Value = -5974

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 25452
$ws.Range("I75").Value = 15500
$ws.Range("J75").Value = 32086.666
$ws.Range("K75").Value = 15500
$ws.Range("L75").Value = 32086.666
$ws.Range("M75").Value = -14564
$ws.Range("N75").Value = -33958.666

$ws.Range("H78").Value = 25452
$ws.Range("I78").Value = 15500
$ws.Range("J78").Value = 32086.666
$ws.Range("K78").Value = 46500
$ws.Range("L78").Value = 96259.99800000001
$ws.Range("M78").Value = -41820
$ws.Range("N78").Value = -105619.998

$ws.Range("H81").Value = 4500
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 4500
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608
